$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 864.4074000000001
$ws.Range("I33").Value = 922.5454999999999
$ws.Range("K33").Value = 922.5454999999999
$ws.Range("M33").Value = -693.5454999999999
$ws.Range("H40").Value = 4148.625
$ws.Range("I40").Value = 3865
$ws.Range("K40").Value = 3865
$ws.Range("M40").Value = -3690
$ws.Range("H53").Value = 490.0625
$ws.Range("J53").Value = 279
$ws.Range("L53").Value = 279
$ws.Range("N53").Value = -1553
$ws.Range("H69").Value = 11069.23
$ws.Range("I69").Value = 6766.778
$ws.Range("K69").Value = 20300.334
$ws.Range("M69").Value = -19426.334
$ws.Range("H70").Value = 1234.1538
$ws.Range("H72").Value = 11069.23
$ws.Range("I72").Value = 6766.778
$ws.Range("K72").Value = 60901.002
$ws.Range("M72").Value = -56533.002
$ws.Range("H73").Value = 1234.1538
$ws.Range("H74").Value = 4987.0625
$ws.Range("I74").Value = 4816.0835
$ws.Range("K74").Value = 4816.0835
$ws.Range("M74").Value = -3880.0835
$ws.Range("H77").Value = 4987.0625
$ws.Range("I77").Value = 4816.0835
$ws.Range("K77").Value = 24080.4175
$ws.Range("M77").Value = -19400.4175
$ws.Range("H86").Value = 6805
$ws.Range("I86").Value = 6696.091
$ws.Range("K86").Value = 6696.091
$ws.Range("M86").Value = -5573.091
$ws.Range("H89").Value = 6805
$ws.Range("I89").Value = 6696.091
$ws.Range("K89").Value = 33480.455
$ws.Range("M89").Value = -27864.455
$ws.Range("H98").Value = 952.84784
$ws.Range("I98").Value = 1071.0769
$ws.Range("K98").Value = 1071.0769
$ws.Range("M98").Value = 426.9231
$ws.Range("H112").Value = 2174.7144
$ws.Range("I112").Value = 1399
$ws.Range("J112").Value = 2267.8
$ws.Range("K112").Value = 4197
$ws.Range("L112").Value = 6803.400000000001
$ws.Range("M112").Value = -3089
$ws.Range("N112").Value = -9019.400000000001
$ws.Range("H116").Value = 5200.3076
$ws.Range("J116").Value = 2603
$ws.Range("L116").Value = 2603
$ws.Range("N116").Value = -9487
$ws.Range("H121").Value = 5061.8
$ws.Range("J121").Value = 4302.25
$ws.Range("L121").Value = 12906.75
$ws.Range("N121").Value = -16400.75
$ws.Range("H122").Value = 952.84784
$ws.Range("I122").Value = 1071.0769
$ws.Range("K122").Value = 3213.2307
$ws.Range("M122").Value = -763.2307000000001
$ws.Range("H131").Value = 4881.7617
$ws.Range("I131").Value = 1432.3125
$ws.Range("K131").Value = 4296.9375
$ws.Range("M131").Value = 743.0625
$ws.Range("H132").Value = 5941.685
$ws.Range("I132").Value = 3249.681
$ws.Range("J132").Value = 24016.572
$ws.Range("K132").Value = 9749.043
$ws.Range("L132").Value = 72049.716
$ws.Range("M132").Value = -7219.043
$ws.Range("N132").Value = -77109.716
$ws.Range("H135").Value = 2050.1738
$ws.Range("I135").Value = 1547.6471
$ws.Range("J135").Value = 3474
$ws.Range("K135").Value = 13928.8239
$ws.Range("L135").Value = 31266
$ws.Range("M135").Value = -11393.8239
$ws.Range("N135").Value = -36336
$ws.Range("H137").Value = 3331.5967
$ws.Range("I137").Value = 3640.625
$ws.Range("J137").Value = 2272.0715
$ws.Range("K137").Value = 10921.875
$ws.Range("L137").Value = 6816.2145
$ws.Range("M137").Value = -8371.875
$ws.Range("N137").Value = -11916.2145
$ws.Range("H141").Value = 14107
$ws.Range("I141").Value = 14525.059
$ws.Range("K141").Value = 43575.177
$ws.Range("M141").Value = -38395.177

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5559.377
$ws.Range("I32").Value = 4949.956
$ws.Range("K32").Value = 4949.956
$ws.Range("M32").Value = -4662.956
$ws.Range("H35").Value = 11666
$ws.Range("I35").Value = 14998
$ws.Range("J35").Value = 10000
$ws.Range("K35").Value = 14998
$ws.Range("L35").Value = 10000
$ws.Range("M35").Value = -14592
$ws.Range("N35").Value = -10812
$ws.Range("H45").Value = 7706.16
$ws.Range("I45").Value = 7662.609
$ws.Range("K45").Value = 7662.609
$ws.Range("M45").Value = -7285.609
$ws.Range("H61").Value = 3092.3098
$ws.Range("I61").Value = 2832.3281
$ws.Range("K61").Value = 2832.3281
$ws.Range("M61").Value = -2620.3281
$ws.Range("H74").Value = 1953.4615
$ws.Range("I74").Value = 1908
$ws.Range("K74").Value = 1908
$ws.Range("M74").Value = -1034
$ws.Range("H77").Value = 1953.4615
$ws.Range("I77").Value = 1908
$ws.Range("K77").Value = 9540
$ws.Range("M77").Value = -5172
$ws.Range("H122").Value = 12928.969
$ws.Range("I122").Value = 2314.0417
$ws.Range("J122").Value = 44773.75
$ws.Range("K122").Value = 6942.125100000001
$ws.Range("L122").Value = 134321.25
$ws.Range("M122").Value = -4492.125100000001
$ws.Range("N122").Value = -139221.25
$ws.Range("H131").Value = 200000
$ws.Range("J131").Value = 200000
$ws.Range("L131").Value = 200000
$ws.Range("N131").Value = -210080
$ws.Range("H136").Value = 3092.3098
$ws.Range("I136").Value = 2832.3281
$ws.Range("K136").Value = 8496.9843
$ws.Range("M136").Value = -5946.9843

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1317.7222
$ws.Range("I20").Value = 1502.5769
$ws.Range("J20").Value = 837.1
$ws.Range("K20").Value = 1502.5769
$ws.Range("L20").Value = 837.1
$ws.Range("M20").Value = -1255.5769
$ws.Range("N20").Value = -1331.1
$ws.Range("H22").Value = 599.5
$ws.Range("I22").Value = 199
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 199
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -26
$ws.Range("N22").Value = -1346
$ws.Range("H80").Value = 4763645.5
$ws.Range("I80").Value = 1456
$ws.Range("J80").Value = 6668521
$ws.Range("K80").Value = 1456
$ws.Range("L80").Value = 6668521
$ws.Range("M80").Value = -458
$ws.Range("N80").Value = -6670517
$ws.Range("H83").Value = 4763645.5
$ws.Range("I83").Value = 1456
$ws.Range("J83").Value = 6668521
$ws.Range("K83").Value = 7280
$ws.Range("L83").Value = 33342605
$ws.Range("M83").Value = -2288
$ws.Range("N83").Value = -33352589
$ws.Range("H99").Value = 2760.8823
$ws.Range("I99").Value = 694.3333
$ws.Range("J99").Value = 5085.75
$ws.Range("K99").Value = 694.3333
$ws.Range("L99").Value = 5085.75
$ws.Range("M99").Value = 803.6667
$ws.Range("N99").Value = -8081.75
$ws.Range("H105").Value = 2244.7222
$ws.Range("I105").Value = 2189.1875
$ws.Range("K105").Value = 2189.1875
$ws.Range("M105").Value = -442.1875
$ws.Range("H134").Value = 3822.4666
$ws.Range("I134").Value = 3819.1162
$ws.Range("K134").Value = 11457.3486
$ws.Range("M134").Value = -8922.348599999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 10322.75
$ws.Range("I22").Value = 19947
$ws.Range("J22").Value = 698.5
$ws.Range("K22").Value = 19947
$ws.Range("L22").Value = 698.5
$ws.Range("M22").Value = -19597
$ws.Range("N22").Value = -1398.5
$ws.Range("H31").Value = 1747.8235
$ws.Range("I31").Value = 1682.125
$ws.Range("J31").Value = 2799
$ws.Range("K31").Value = 1682.125
$ws.Range("L31").Value = 2799
$ws.Range("M31").Value = -1387.125
$ws.Range("N31").Value = -3389
$ws.Range("H34").Value = 1747.8235
$ws.Range("I34").Value = 1682.125
$ws.Range("J34").Value = 2799
$ws.Range("K34").Value = 1682.125
$ws.Range("L34").Value = 2799
$ws.Range("M34").Value = -1480.125
$ws.Range("N34").Value = -3203
$ws.Range("H35").Value = 1133.125
$ws.Range("I35").Value = 580.7143
$ws.Range("J35").Value = 5000
$ws.Range("K35").Value = 580.7143
$ws.Range("L35").Value = 5000
$ws.Range("M35").Value = -286.7143
$ws.Range("N35").Value = -5588
$ws.Range("H63").Value = 52499.25
$ws.Range("J63").Value = 52499.25
$ws.Range("L63").Value = 52499.25
$ws.Range("N63").Value = -53871.25
$ws.Range("H66").Value = 52499.25
$ws.Range("J66").Value = 52499.25
$ws.Range("L66").Value = 157497.75
$ws.Range("N66").Value = -164361.75
$ws.Range("H68").Value = 47340
$ws.Range("J68").Value = 46729.25
$ws.Range("L68").Value = 46729.25
$ws.Range("N68").Value = -48227.25
$ws.Range("H71").Value = 47340
$ws.Range("J71").Value = 46729.25
$ws.Range("L71").Value = 140187.75
$ws.Range("N71").Value = -147675.75
$ws.Range("H86").Value = 6063659.5
$ws.Range("I86").Value = 9526637
$ws.Range("J86").Value = 3449.25
$ws.Range("K86").Value = 9526637
$ws.Range("L86").Value = 3449.25
$ws.Range("M86").Value = -9525514
$ws.Range("N86").Value = -5695.25
$ws.Range("H89").Value = 6063659.5
$ws.Range("I89").Value = 9526637
$ws.Range("J89").Value = 3449.25
$ws.Range("K89").Value = 47633185
$ws.Range("L89").Value = 17246.25
$ws.Range("M89").Value = -47627569
$ws.Range("N89").Value = -28478.25
$ws.Range("H105").Value = 1926.2727
$ws.Range("I105").Value = 999.2
$ws.Range("J105").Value = 2698.8333
$ws.Range("K105").Value = 999.2
$ws.Range("L105").Value = 2698.8333
$ws.Range("M105").Value = 747.8
$ws.Range("N105").Value = -6192.8333
$ws.Range("H132").Value = 4202.6665
$ws.Range("I132").Value = 2194.2563
$ws.Range("J132").Value = 12905.777
$ws.Range("K132").Value = 6582.7689
$ws.Range("L132").Value = 38717.331
$ws.Range("M132").Value = -4052.7689
$ws.Range("N132").Value = -43777.331
$ws.Range("H134").Value = 2177.3022
$ws.Range("I134").Value = 2110.0952
$ws.Range("K134").Value = 6330.285600000001
$ws.Range("M134").Value = -3795.285600000001
$ws.Range("H135").Value = 79610.14
$ws.Range("J135").Value = 79610.14
$ws.Range("L135").Value = 79610.14
$ws.Range("N135").Value = -89750.14
$ws.Range("H138").Value = 76797.664
$ws.Range("J138").Value = 76797.664
$ws.Range("L138").Value = 76797.664
$ws.Range("N138").Value = -87077.664
$ws.Range("H141").Value = 121846.27
$ws.Range("J141").Value = 128106.16
$ws.Range("L141").Value = 128106.16
$ws.Range("N141").Value = -138466.16

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 30686.273
$ws.Range("I4").Value = 140.34616
$ws.Range("J4").Value = 144142.58
$ws.Range("K4").Value = 421.03848
$ws.Range("L4").Value = 432427.74
$ws.Range("M4").Value = -309.03848
$ws.Range("N4").Value = -432651.74
$ws.Range("H5").Value = 785.2857
$ws.Range("I5").Value = 785.2857
$ws.Range("K5").Value = 2355.8571
$ws.Range("M5").Value = -2243.8571
$ws.Range("H14").Value = 253
$ws.Range("I14").Value = 253
$ws.Range("K14").Value = 759
$ws.Range("M14").Value = -586
$ws.Range("H56").Value = 11601.4375
$ws.Range("I56").Value = 11601.4375
$ws.Range("K56").Value = 11601.4375
$ws.Range("M56").Value = -11071.4375
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").ClearContents()
$ws.Range("H131").Value = 818806.3
$ws.Range("I131").Value = 2942282.2
$ws.Range("K131").Value = 8826846.600000001
$ws.Range("M131").Value = -8821806.600000001
$ws.Range("H135").Value = 785.2857
$ws.Range("I135").Value = 785.2857
$ws.Range("K135").Value = 7067.571300000001
$ws.Range("M135").Value = -4532.571300000001
$ws.Range("H137").Value = 1846.2307
$ws.Range("I137").Value = 1846.2307
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 5538.6921
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -438.6921000000002
$ws.Range("N137").ClearContents()
$ws.Range("H140").Value = 1677.5
$ws.Range("I140").Value = 1677.5
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 5032.5
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = 147.5
$ws.Range("N140").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("H97").Value = 180.0625
$ws.Range("I97").Value = 165.6
$ws.Range("J97").Value = 186.63637
$ws.Range("K97").Value = 165.6
$ws.Range("L97").Value = 186.63637
$ws.Range("M97").Value = 330.4
$ws.Range("N97").Value = -1178.63637
$ws.Range("H113").Value = 2358
$ws.Range("I113").Value = 2155.7222
$ws.Range("K113").Value = 2155.7222
$ws.Range("M113").Value = 14.27779999999984
$ws.Range("H122").Value = 3351.4119
$ws.Range("I122").Value = 2750.5833
$ws.Range("J122").Value = 4793.4
$ws.Range("K122").Value = 8251.749899999999
$ws.Range("L122").Value = 14380.2
$ws.Range("M122").Value = -5801.749899999999
$ws.Range("N122").Value = -19280.2
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 8334504.5
$ws.Range("I16").Value = 10417748
$ws.Range("J16").Value = 1529.3334
$ws.Range("K16").Value = 10417748
$ws.Range("L16").Value = 1529.3334
$ws.Range("M16").Value = -10417578
$ws.Range("N16").Value = -1869.3334
$ws.Range("H23").Value = 33000
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H46").Value = 3677.7144
$ws.Range("I46").Value = 1033
$ws.Range("J46").Value = 4399
$ws.Range("K46").Value = 1033
$ws.Range("L46").Value = 4399
$ws.Range("M46").Value = -845
$ws.Range("N46").Value = -4775
$ws.Range("H122").Value = 3534.75
$ws.Range("I122").Value = 3534.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10604.25
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8154.25
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 40299.55
$ws.Range("I132").Value = 48284.41
$ws.Range("J132").Value = 15204.286
$ws.Range("K132").Value = 144853.23
$ws.Range("L132").Value = 45612.858
$ws.Range("M132").Value = -142323.23
$ws.Range("N132").Value = -50672.858
$ws.Range("H136").Value = 11304.9375
$ws.Range("J136").Value = 5790.1665
$ws.Range("L136").Value = 17370.4995
$ws.Range("N136").Value = -22470.4995

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 727.7
$ws.Range("I100").Value = 727.7
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1455.4
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -914.4000000000001
$ws.Range("N100").ClearContents()
$ws.Range("H132").Value = 9724.625
$ws.Range("I132").Value = 10966
$ws.Range("J132").Value = 8979.799999999999
$ws.Range("K132").Value = 32898
$ws.Range("L132").Value = 26939.4
$ws.Range("M132").Value = -30368
$ws.Range("N132").Value = -31999.4
$ws.Range("H136").Value = 2419
$ws.Range("I136").Value = 2840.842
$ws.Range("J136").Value = 1083.1666
$ws.Range("K136").Value = 8522.526
$ws.Range("L136").Value = 3249.4998
$ws.Range("M136").Value = -5972.526
$ws.Range("N136").Value = -8349.4998
$ws.Range("H137").Value = 147500
$ws.Range("J137").Value = 147500
$ws.Range("L137").Value = 147500
$ws.Range("N137").Value = -157700
